$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "12:20-15:45"

# Enter the date as a formula that evaluates to text, then convert it to a
# plain value in-place via copy/paste-special. This avoids Excel's
# autodetection turning the literal "01.11.2023" into a real date serial
# (and avoids leaving behind an unused number-format style).
$ws.Range("A4").Formula = '="01.11.2023"'
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163)

$ws.Range("C4").Value = "Getting started on coding"

$ws.Range("C4").Select()
